$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 182.05556
$ws.Range("I33").Value = 104.75
$ws.Range("K33").Value = 104.75
$ws.Range("M33").Value = 124.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1508.55
$ws.Range("I40").Value = 1097.9231
$ws.Range("J40").Value = 2271.1428
$ws.Range("K40").Value = 1097.9231
$ws.Range("L40").Value = 2271.1428
$ws.Range("M40").Value = -922.9231
$ws.Range("N40").Value = -2621.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 467.4375
$ws.Range("I53").Value = 645.5
$ws.Range("J53").Value = 170.66667
$ws.Range("K53").Value = 645.5
$ws.Range("L53").Value = 170.66667
$ws.Range("M53").Value = -8.5
$ws.Range("N53").Value = -1444.66667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 159.0625
$ws.Range("J55").Value = 198
$ws.Range("L55").Value = 198
$ws.Range("N55").Value = -626

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 6000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1969.85
$ws.Range("I70").Value = 1599
$ws.Range("J70").Value = 2011.0555
$ws.Range("K70").Value = 4797
$ws.Range("L70").Value = 6033.166499999999
$ws.Range("M70").Value = -4527
$ws.Range("N70").Value = -6573.166499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 6000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1969.85
$ws.Range("I73").Value = 1599
$ws.Range("J73").Value = 2011.0555
$ws.Range("K73").Value = 4797
$ws.Range("L73").Value = 6033.166499999999
$ws.Range("M73").Value = -3861
$ws.Range("N73").Value = -7905.166499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 637.8
$ws.Range("I107").Value = 577.25
$ws.Range("J107").Value = 880
$ws.Range("K107").Value = 577.25
$ws.Range("L107").Value = 880
$ws.Range("M107").Value = 1342.75
$ws.Range("N107").Value = -4720

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1999.091
$ws.Range("J138").Value = 2224.25
$ws.Range("L138").Value = 6672.75
$ws.Range("N138").Value = -16952.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1677.9231
$ws.Range("I97").Value = 1601.125
$ws.Range("K97").Value = 1601.125
$ws.Range("M97").Value = -1105.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 983.6
$ws.Range("I102").Value = 976.75
$ws.Range("K102").Value = 976.75
$ws.Range("M102").Value = 645.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4473.3335
$ws.Range("I86").Value = 752
$ws.Range("K86").Value = 752
$ws.Range("M86").Value = 371

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4473.3335
$ws.Range("I89").Value = 752
$ws.Range("K89").Value = 3760
$ws.Range("M89").Value = 1856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1535.125
$ws.Range("I99").Value = 1576.8667
$ws.Range("K99").Value = 1576.8667
$ws.Range("M99").Value = -78.86670000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1674.4706
$ws.Range("I7").Value = 943.53845
$ws.Range("J7").Value = 4050
$ws.Range("K7").Value = 943.53845
$ws.Range("L7").Value = 4050
$ws.Range("M7").Value = -830.53845
$ws.Range("N7").Value = -4276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1917.7693
$ws.Range("I22").Value = 628.875
$ws.Range("J22").Value = 3980
$ws.Range("K22").Value = 628.875
$ws.Range("L22").Value = 3980
$ws.Range("M22").Value = -278.875
$ws.Range("N22").Value = -4680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3359.8
$ws.Range("I62").Value = 2949.75
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2949.75
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2325.75
$ws.Range("N62").Value = -6248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3359.8
$ws.Range("I65").Value = 2949.75
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 14748.75
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -11628.75
$ws.Range("N65").Value = -31240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 684.0833
$ws.Range("I5").Value = 822.25
$ws.Range("J5").Value = 615
$ws.Range("K5").Value = 2466.75
$ws.Range("L5").Value = 1845
$ws.Range("M5").Value = -2354.75
$ws.Range("N5").Value = -2069

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 905
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 999.5
$ws.Range("J39").Value = 999.5
$ws.Range("L39").Value = 2998.5
$ws.Range("N39").Value = -3586.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 8520
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 10550
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 31650
$ws.Range("M46").Value = -1109
$ws.Range("N46").Value = -31832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1995
$ws.Range("I70").Value = 1995
$ws.Range("K70").Value = 5985
$ws.Range("M70").Value = -5670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1995
$ws.Range("I73").Value = 1995
$ws.Range("K73").Value = 5985
$ws.Range("M73").Value = -4893

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2500
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7500
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = -6377
$ws.Range("N81").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2500
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 22500
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = -16884
$ws.Range("N84").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 7600
$ws.Range("J93").Value = 15000
$ws.Range("L93").Value = 45000
$ws.Range("N93").Value = -48744

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 6753.6665
$ws.Range("J104").Value = 6753.6665
$ws.Range("L104").Value = 20260.9995
$ws.Range("N104").Value = -25502.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 2000
$ws.Range("K115").Value = 6000
$ws.Range("M115").Value = -4825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 962.9474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 684.0833
$ws.Range("I135").Value = 822.25
$ws.Range("J135").Value = 615
$ws.Range("K135").Value = 7400.25
$ws.Range("L135").Value = 5535
$ws.Range("M135").Value = -4865.25
$ws.Range("N135").Value = -10605

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 700
$ws.Range("I138").Value = 700
$ws.Range("K138").Value = 2100
$ws.Range("M138").Value = 3040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6538.75
$ws.Range("I139").Value = 5385
$ws.Range("K139").Value = 16155
$ws.Range("M139").Value = -11015

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1265
$ws.Range("I140").Value = 1265
$ws.Range("K140").Value = 3795
$ws.Range("M140").Value = 1385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3424
$ws.Range("I141").Value = 3424
$ws.Range("K141").Value = 10272
$ws.Range("M141").Value = -5092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 6669.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4893.3335
$ws.Range("I113").Value = 4893.3335
$ws.Range("K113").Value = 4893.3335
$ws.Range("M113").Value = -2723.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7379.4
$ws.Range("I126").Value = 6224.25
$ws.Range("K126").Value = 18672.75
$ws.Range("M126").Value = -16202.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 6750
$ws.Range("J6").Value = 6750
$ws.Range("L6").Value = 6750
$ws.Range("N6").Value = -6980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 656.5
$ws.Range("I81").Value = 609.75
$ws.Range("J81").Value = 750
$ws.Range("K81").Value = 1219.5
$ws.Range("L81").Value = 1500
$ws.Range("M81").Value = -158.5
$ws.Range("N81").Value = -3622

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 656.5
$ws.Range("I84").Value = 609.75
$ws.Range("J84").Value = 750
$ws.Range("K84").Value = 6097.5
$ws.Range("L84").Value = 7500
$ws.Range("M84").Value = -793.5
$ws.Range("N84").Value = -18108

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2911
$ws.Range("J132").Value = 3062.375
$ws.Range("L132").Value = 9187.125
$ws.Range("N132").Value = -14247.125
